$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header rows 19-20 (merged header block) + data rows 21-41
# ---------------------------------------------------------------------------

# Row 19 - top header
$ws.Range("A19").Value = "Name"
$ws.Range("B19").Value = "accurate circuit"
$ws.Range("F19").Value = "approximate circuit"

# Row 20 - column headers
$ws.Range("B20").Value = "I/O"
$ws.Range("C20").Value = "\#nodes"
$ws.Range("D20").Value = "Original area"
$ws.Range("E20").Value = "Original delay"
$ws.Range("F20").Value = "Area ratio"
$ws.Range("G20").Value = "Delay ratio"
$ws.Range("H20").Value = "Runtime/s"

# Row 21
$ws.Range("A21").Value = "adder"
$ws.Range("B21").Value = "256/129"
$ws.Range("C21").Value = 1117
$ws.Range("D21").Value = 2594
$ws.Range("E21").Value = 303.8

# Row 22
$ws.Range("A22").Value = "arbiter"
$ws.Range("B22").Value = "256/129"
$ws.Range("C22").Value = 857
$ws.Range("D22").Value = 1802
$ws.Range("E22").Value = 242.6

# Row 23
$ws.Range("A23").Value = "bar"
$ws.Range("B23").Value = "135/128"
$ws.Range("C23").Value = 1927
$ws.Range("D23").Value = 5383
$ws.Range("E23").Value = 51.6

# Row 24
$ws.Range("A24").Value = "cavlc"
$ws.Range("B24").Value = "10/11"
$ws.Range("C24").Value = 441
$ws.Range("D24").Value = 1093
$ws.Range("E24").Value = 24.2

# Row 25
$ws.Range("A25").Value = "ctrl"
$ws.Range("B25").Value = "7/26"
$ws.Range("C25").Value = 87
$ws.Range("D25").Value = 195
$ws.Range("E25").Value = 12.7

# Row 26
$ws.Range("A26").Value = "dec"
$ws.Range("B26").Value = "8/256"
$ws.Range("C26").Value = 433
$ws.Range("D26").Value = 995
$ws.Range("E26").Value = 29

# Row 27 (D27 holds the literal formula =47469)
$ws.Range("A27").Value = "div"
$ws.Range("B27").Value = "128/128"
$ws.Range("C27").Value = 17710
$ws.Range("D27").Formula = "=47469"
$ws.Range("E27").Value = 5533.8

# Row 28
$ws.Range("A28").Value = "hyp"
$ws.Range("B28").Value = "256/128"
$ws.Range("C28").Value = 278279
$ws.Range("D28").Value = 687703
$ws.Range("E28").Value = 17121.19

# Row 29
$ws.Range("A29").Value = "i2c"
$ws.Range("B29").Value = "147/142"
$ws.Range("C29").Value = 611
$ws.Range("D29").Value = 1428
$ws.Range("E29").Value = 31.2

# Row 30
$ws.Range("A30").Value = "int2float"
$ws.Range("B30").Value = "11/7"
$ws.Range("C30").Value = 147
$ws.Range("D30").Value = 341
$ws.Range("E30").Value = 19.5

# Row 31
$ws.Range("A31").Value = "log"
$ws.Range("B31").Value = "32/32"
$ws.Range("C31").Value = 27468
$ws.Range("D31").Value = 69688
$ws.Range("E31").Value = 651.4

# Row 32
$ws.Range("A32").Value = "max"
$ws.Range("B32").Value = "512/130"
$ws.Range("C32").Value = 2163
$ws.Range("D32").Value = 5041
$ws.Range("E32").Value = 466.9

# Row 33
$ws.Range("A33").Value = "mem"
$ws.Range("B33").Value = "1204/1231"
$ws.Range("C33").Value = 6205
$ws.Range("D33").Value = 14671
$ws.Range("E33").Value = 101.9

# Row 34
$ws.Range("A34").Value = "multiplier"
$ws.Range("B34").Value = "128/128"
$ws.Range("C34").Value = 20260
$ws.Range("D34").Value = 54205
$ws.Range("E34").Value = 419.5

# Row 35
$ws.Range("A35").Value = "priority"
$ws.Range("B35").Value = "128/8"
$ws.Range("C35").Value = 351
$ws.Range("D35").Value = 741
$ws.Range("E35").Value = 126.8

# Row 36
$ws.Range("A36").Value = "router"
$ws.Range("B36").Value = "60/30"
$ws.Range("C36").Value = 109
$ws.Range("D36").Value = 186
$ws.Range("E36").Value = 13.7

# Row 37
$ws.Range("A37").Value = "sin"
$ws.Range("B37").Value = "24/25"
$ws.Range("C37").Value = 5534
$ws.Range("D37").Value = 13552
$ws.Range("E37").Value = 272.9

# Row 38
$ws.Range("A38").Value = "sqrt"
$ws.Range("B38").Value = "128/64"
$ws.Range("C38").Value = 16584
$ws.Range("D38").Value = 43859
$ws.Range("E38").Value = 7304

# Row 39
$ws.Range("A39").Value = "square"
$ws.Range("B39").Value = "64/128"
$ws.Range("C39").Value = 14967
$ws.Range("D39").Value = 37672
$ws.Range("E39").Value = 355.5

# Row 40
$ws.Range("A40").Value = "voter"
$ws.Range("B40").Value = "1001/1"
$ws.Range("C40").Value = 14112
$ws.Range("D40").Value = 33683
$ws.Range("E40").Value = 95.8

# Row 41
$ws.Range("A41").Value = "average"

# ---------------------------------------------------------------------------
# 2. Formatting: old leftover rows 24-31 (columns D/E) used the plain
#    center-text style; clear F-H for rows where nothing was written so that
#    no stale style lingers and every touched cell picks up the uniform
#    center alignment applied below.
# ---------------------------------------------------------------------------

# Make sure F:H exist (empty but styled) for rows 19, 21-41
$ws.Range("F21:H41").Value = ""

# ---------------------------------------------------------------------------
# 3. Apply the center/center alignment to the whole new block A19:H41
#    (this is cellXfs style index 4 in the target workbook)
# ---------------------------------------------------------------------------
$block = $ws.Range("A19:H41")
$block.HorizontalAlignment = -4108
$block.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. The "I/O" fraction-like text column (B) additionally uses the text
#    number format "@" (cellXfs style index 5), except for the header cells
#    B19/B20 and the B28 "hyp" row, which stay on the plain style (4).
# ---------------------------------------------------------------------------
$ws.Range("B21:B27").NumberFormat = "@"
$ws.Range("B29:B40").NumberFormat = "@"
$ws.Range("B21:B27").HorizontalAlignment = -4108
$ws.Range("B21:B27").VerticalAlignment = -4108
$ws.Range("B29:B40").HorizontalAlignment = -4108
$ws.Range("B29:B40").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Merge the header cells
# ---------------------------------------------------------------------------
$ws.Range("B19:E19").Merge()
$ws.Range("F19:H19").Merge()
$ws.Range("A19:A20").Merge()

# ---------------------------------------------------------------------------
# 6. Column widths for C, D, E
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 11.97
$ws.Columns.Item(4).ColumnWidth = 10.37
$ws.Columns.Item(5).ColumnWidth = 11.37

# ---------------------------------------------------------------------------
# 7. View state: scrolled to A9, active cell E21
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E21").Select()
